$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value2 = 0.07666397707932049
$ws.Range("D2").Value2 = 0.1304900858569109
$ws.Range("E2").Value2 = 0.07707643088466121
$ws.Range("F2").Value2 = 2.379194894897779
$ws.Range("G2").Value2 = 0.00256018347049029
$ws.Range("I2").Value2 = 2.190033717674908
$ws.Range("K2").Value2 = 2.397405336578231
$ws.Range("L2").Value2 = 0.1216957980070674
$ws.Range("M2").Value2 = 0.6167442470982891
$ws.Range("N2").Value2 = 2.060436144321823

$ws.Range("C3").Value2 = 0.07609603802299603
$ws.Range("D3").Value2 = 0.1318212558796983
$ws.Range("E3").Value2 = 0.07700476154855274
$ws.Range("F3").Value2 = 2.361058283553419
$ws.Range("G3").Value2 = 0.002565446804889912
$ws.Range("I3").Value2 = 2.179414616119942
$ws.Range("K3").Value2 = 2.262280513432131
$ws.Range("L3").Value2 = 0.1219348201225152
$ws.Range("M3").Value2 = 0.5922771613383091
$ws.Range("N3").Value2 = 2.084487371412415

$ws.Range("C4").Value2 = 0.07577229705513844
$ws.Range("D4").Value2 = 0.1326775970894847
$ws.Range("E4").Value2 = 0.0769907625790669
$ws.Range("F4").Value2 = 2.351385231975655
$ws.Range("G4").Value2 = 0.002568848142911478
$ws.Range("I4").Value2 = 2.174171333713943
$ws.Range("K4").Value2 = 2.180622279277827
$ws.Range("L4").Value2 = 0.1221326516306682
$ws.Range("M4").Value2 = 0.5775979327944469
$ws.Range("N4").Value2 = 2.099982966035345

$ws.Range("C5").Value2 = 0.07564667176072248
$ws.Range("D5").Value2 = 0.1330363611740228
$ws.Range("E5").Value2 = 0.07699260667528307
$ws.Range("F5").Value2 = 2.347809934170527
$ws.Range("G5").Value2 = 0.002570277016704092
$ws.Range("I5").Value2 = 2.172354746305743
$ws.Range("K5").Value2 = 2.147673813637368
$ws.Range("L5").Value2 = 0.1222261031390879
$ws.Range("M5").Value2 = 0.5717022369385987
$ws.Range("N5").Value2 = 2.106480344901101

$ws.Range("C6").Value2 = 0.07562619302445484
$ws.Range("D6").Value2 = 0.1330965252274838
$ws.Range("E6").Value2 = 0.07699336892258302
$ws.Range("F6").Value2 = 2.347238359396755
$ws.Range("G6").Value2 = 0.002570516869438923
$ws.Range("I6").Value2 = 2.172072407435365
$ws.Range("K6").Value2 = 2.1422225173161
$ws.Range("L6").Value2 = 0.1222423955624805
$ws.Range("M6").Value2 = 0.5707284620650412
$ws.Range("N6").Value2 = 2.107570262348315

$ws.Range("C7").Value2 = 0.0757705772827677
$ws.Range("D7").Value2 = 0.132682395847139
$ws.Range("E7").Value2 = 0.07699075687944479
$ws.Range("F7").Value2 = 2.351335531862077
$ws.Range("G7").Value2 = 0.00256886723981083
$ws.Range("I7").Value2 = 2.174145539862437
$ws.Range("K7").Value2 = 2.18017659857486
$ws.Range("L7").Value2 = 0.1221338599989998
$ws.Range("M7").Value2 = 0.57751807259217
$ws.Range("N7").Value2 = 2.100069852197358

$ws.Range("C8").Value2 = 0.07646298035036381
$ws.Range("D8").Value2 = 0.1309409745263217
$ws.Range("E8").Value2 = 0.07704549330775912
$ws.Range("F8").Value2 = 2.372636962133981
$ws.Range("G8").Value2 = 0.00256196314987807
$ws.Range("I8").Value2 = 2.186106623958437
$ws.Range("K8").Value2 = 2.350541833068291
$ws.Range("L8").Value2 = 0.1217676075421643
$ws.Range("M8").Value2 = 0.6082365330170418
$ws.Range("N8").Value2 = 2.068577738451925

$ws.Range("C9").Value2 = 0.07801815990123373
$ws.Range("D9").Value2 = 0.1278357976846127
$ws.Range("E9").Value2 = 0.07739086237637771
$ws.Range("F9").Value2 = 2.426081339935834
$ws.Range("G9").Value2 = 0.002549763495880404
$ws.Range("I9").Value2 = 2.219743040017391
$ws.Range("K9").Value2 = 2.695089803630196
$ws.Range("L9").Value2 = 0.1214551154872332
$ws.Range("M9").Value2 = 0.6712156724132257
$ws.Range("N9").Value2 = 2.012609857049913

$ws.Range("C10").Value2 = 0.07928022775045207
$ws.Range("D10").Value2 = 0.1257436242040324
$ws.Range("E10").Value2 = 0.07778978510310175
$ws.Range("F10").Value2 = 2.472560946395703
$ws.Range("G10").Value2 = 0.00254160747538694
$ws.Range("I10").Value2 = 2.250736873597532
$ws.Range("K10").Value2 = 2.954743712645154
$ws.Range("L10").Value2 = 0.1214737303910169
$ws.Range("M10").Value2 = 0.7191818108004355
$ws.Range("N10").Value2 = 1.975033076299418

$ws.Range("C11").Value2 = 0.07988015856168573
$ws.Range("D11").Value2 = 0.1248330100378396
$ws.Range("E11").Value2 = 0.07800281982671819
$ws.Range("F11").Value2 = 2.495294183631131
$ws.Range("G11").Value2 = 0.00253807033475818
$ws.Range("I11").Value2 = 2.266217430370176
$ws.Range("K11").Value2 = 3.074313474267058
$ws.Range("L11").Value2 = 0.1215362982367481
$ws.Range("M11").Value2 = 0.7413768248265882
$ws.Range("N11").Value2 = 1.95871086219552

$ws.Range("C12").Value2 = 0.08011103145635445
$ws.Range("D12").Value2 = 0.124494109094158
$ws.Range("E12").Value2 = 0.07808802981015006
$ws.Range("F12").Value2 = 2.504132896327064
$ws.Range("G12").Value2 = 0.002536755646681408
$ws.Range("I12").Value2 = 2.272279421085784
$ws.Range("K12").Value2 = 3.119802239953231
$ws.Range("L12").Value2 = 0.1215677857119672
$ws.Range("M12").Value2 = 0.7498357870869086
$ws.Range("N12").Value2 = 1.952641352760029

$ws.Range("C13").Value2 = 0.08006114494203587
$ws.Range("D13").Value2 = 0.1245668335421719
$ws.Range("E13").Value2 = 0.07806947647283735
$ws.Range("F13").Value2 = 2.502219061078762
$ws.Range("G13").Value2 = 0.002537037689822558
$ws.Range("I13").Value2 = 2.27096495551838
$ws.Range("K13").Value2 = 3.109996042997011
$ws.Range("L13").Value2 = 0.1215606574202113
$ws.Range("M13").Value2 = 0.748011583585992
$ws.Range("N13").Value2 = 1.953943571091937

$ws.Range("C14").Value2 = 0.07989907872695312
$ws.Range("D14").Value2 = 0.1248050095415589
$ws.Range("E14").Value2 = 0.07800973914860165
$ws.Range("F14").Value2 = 2.496016728633919
$ws.Range("G14").Value2 = 0.002537961679499718
$ws.Range("I14").Value2 = 2.266712141910361
$ws.Range("K14").Value2 = 3.078051639006731
$ws.Range("L14").Value2 = 0.1215387324321711
$ws.Range("M14").Value2 = 0.742071661065097
$ws.Range("N14").Value2 = 1.95820928590549

$ws.Range("C15").Value2 = 0.07980028874712275
$ws.Range("D15").Value2 = 0.124951671789324
$ws.Range("E15").Value2 = 0.07797373933653873
$ws.Range("F15").Value2 = 2.492247637545006
$ws.Range("G15").Value2 = 0.002538530868834101
$ws.Range("I15").Value2 = 2.264133232051108
$ws.Range("K15").Value2 = 3.058512201074109
$ws.Range("L15").Value2 = 0.1215263182548796
$ws.Range("M15").Value2 = 0.7384403526319261
$ws.Range("N15").Value2 = 1.960836673509377

$ws.Range("C16").Value2 = 0.07924153856148308
$ws.Range("D16").Value2 = 0.1258039638106201
$ws.Range("E16").Value2 = 0.07777649787077578
$ws.Range("F16").Value2 = 2.471107393048939
$ws.Range("G16").Value2 = 0.002541842106630148
$ws.Range("I16").Value2 = 2.249753081188842
$ws.Range("K16").Value2 = 2.946958891464647
$ws.Range("L16").Value2 = 0.1214707313232068
$ws.Range("M16").Value2 = 0.7177388828698383
$ws.Range("N16").Value2 = 1.976115332651222

$ws.Range("C17").Value2 = 0.07890536006365068
$ws.Range("D17").Value2 = 0.1263373612049818
$ws.Range("E17").Value2 = 0.07766358026524678
$ws.Range("F17").Value2 = 2.458546739582829
$ws.Range("G17").Value2 = 0.002543917674094889
$ws.Range("I17").Value2 = 2.241285964738708
$ws.Range("K17").Value2 = 2.878897474742814
$ws.Range("L17").Value2 = 0.1214504972407653
$ws.Range("M17").Value2 = 0.7051354320100671
$ws.Range("N17").Value2 = 1.985686212794548

$ws.Range("C18").Value2 = 0.07871442955213581
$ws.Range("D18").Value2 = 0.1266480273230268
$ws.Range("E18").Value2 = 0.07760160388313153
$ws.Range("F18").Value2 = 2.451471693412884
$ws.Range("G18").Value2 = 0.002545127784146519
$ws.Range("I18").Value2 = 2.236545844943549
$ws.Range("K18").Value2 = 2.839886914064778
$ws.Range("L18").Value2 = 0.1214439503419413
$ws.Range("M18").Value2 = 0.6979215382557697
$ws.Range("N18").Value2 = 1.991263737230355

$ws.Range("C19").Value2 = 0.07865020179294646
$ws.Range("D19").Value2 = 0.1267538777682464
$ws.Range("E19").Value2 = 0.0775811300044893
$ws.Range("F19").Value2 = 2.449101835237641
$ws.Range("G19").Value2 = 0.002545540310393034
$ws.Range("I19").Value2 = 2.234963202865714
$ws.Range("K19").Value2 = 2.826702018595711
$ws.Range("L19").Value2 = 0.1214426076489978
$ws.Range("M19").Value2 = 0.6954850886856434
$ws.Range("N19").Value2 = 1.993164649603642

$ws.Range("C20").Value2 = 0.07894089543489713
$ws.Range("D20").Value2 = 0.1262801795002098
$ws.Range("E20").Value2 = 0.07767529307023402
$ws.Range("F20").Value2 = 2.45986835762217
$ws.Range("G20").Value2 = 0.00254369504065184
$ws.Range("I20").Value2 = 2.242173846575568
$ws.Range("K20").Value2 = 2.88612858388143
$ws.Range("L20").Value2 = 0.1214521241665274
$ws.Range("M20").Value2 = 0.7064734378173512
$ws.Range("N20").Value2 = 1.984659858844115

$ws.Range("C21").Value2 = 0.07994658145219091
$ws.Range("D21").Value2 = 0.1247348904469856
$ws.Range("E21").Value2 = 0.07802716228017914
$ws.Range("F21").Value2 = 2.497832245744803
$ws.Range("G21").Value2 = 0.002537689610670246
$ws.Range("I21").Value2 = 2.267955863023616
$ws.Range("K21").Value2 = 3.087428765258267
$ws.Range("L21").Value2 = 0.1215449606734893
$ws.Range("M21").Value2 = 0.7438148873289094
$ws.Range("N21").Value2 = 1.956953316613141

$ws.Range("C22").Value2 = 0.0806253672298709
$ws.Range("D22").Value2 = 0.1237595201541009
$ws.Range("E22").Value2 = 0.07828358224842802
$ws.Range("F22").Value2 = 2.523985955662084
$ws.Range("G22").Value2 = 0.002533908914350454
$ws.Range("I22").Value2 = 2.28597129587159
$ws.Range("K22").Value2 = 3.220216548658982
$ws.Range("L22").Value2 = 0.121651073659983
$ws.Range("M22").Value2 = 0.7685356891852848
$ws.Range("N22").Value2 = 1.939494764129428

$ws.Range("C23").Value2 = 0.08026112441821454
$ws.Range("D23").Value2 = 0.1242769257786343
$ws.Range("E23").Value2 = 0.07814430564933161
$ws.Range("F23").Value2 = 2.509903899750356
$ws.Range("G23").Value2 = 0.002535913594252543
$ws.Range("I23").Value2 = 2.27624908187272
$ws.Range("K23").Value2 = 3.149232474962957
$ws.Range("L23").Value2 = 0.1215902764320163
$ws.Range("M23").Value2 = 0.7553127253750631
$ws.Range("N23").Value2 = 1.948753180255824

$ws.Range("C24").Value2 = 0.07892482259399713
$ws.Range("D24").Value2 = 0.1263060188634491
$ws.Range("E24").Value2 = 0.07766998854910767
$ws.Range("F24").Value2 = 2.45927039871755
$ws.Range("G24").Value2 = 0.002543795640807443
$ws.Range("I24").Value2 = 2.241772037412218
$ws.Range("K24").Value2 = 2.882859029348197
$ws.Range("L24").Value2 = 0.1214513727931035
$ws.Range("M24").Value2 = 0.7058684257453649
$ws.Range("N24").Value2 = 1.985123639581538

$ws.Range("C25").Value2 = 0.07757640127083221
$ws.Range("D25").Value2 = 0.1286426337855033
$ws.Range("E25").Value2 = 0.0772719419448471
$ws.Range("F25").Value2 = 2.410363538931477
$ws.Range("G25").Value2 = 0.002552921421625913
$ws.Range("I25").Value2 = 2.209545922327592
$ws.Range("K25").Value2 = 2.600746453748968
$ws.Range("L25").Value2 = 0.121496129822912
$ws.Range("M25").Value2 = 0.6538822468937653
$ws.Range("N25").Value2 = 2.027129558072446
